# Update on 20210731 画中人
# Replace English double quotes around specific quoted phrases with single
# quotes in the act5d0 "first enter" story lines (Capone / Gambino / Emperor).
#
# Rather than hard-coding row numbers, scan every used cell on the active
# sheet and perform the same four straight-quote -> curly/single-quote style
# swaps the diff describes, only on cells whose text actually contains the
# old (double-quoted) phrasing. This makes the edit resilient to row/column
# shifts while still touching exactly the four cells described in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
$startRow = $used.Row
$startCol = $used.Column

$replacements = @(
    @{ Old = '[name="Capone"]   And besides, here in Lungmen, you and me are equals, "boss."'; `
       New = "[name=`"Capone`"]   And besides, here in Lungmen, you and me are equals, 'boss.'" },
    @{ Old = '[name="Capone"]   Wei Yenwu set some rules for this city. As long as we don''t break ''em, we can run our rackets here just like before. Lungmen needs a little "legitimate business."'; `
       New = "[name=`"Capone`"]   Wei Yenwu set some rules for this city. As long as we don't break 'em, we can run our rackets here just like before. Lungmen needs a little 'legitimate business.'" },
    @{ Old = '[name="Gambino"]   Yo, Capone, do me a favor and torch this tacky joint. Call it "urban beautification." Wei Yenwu''ll thank us for tearin'' down a tastelessly decorated rathole.'; `
       New = "[name=`"Gambino`"]   Yo, Capone, do me a favor and torch this tacky joint. Call it 'urban beautification.' Wei Yenwu'll thank us for tearin' down a tastelessly decorated rathole." },
    @{ Old = '[name="Emperor"]   Who you callin'' "tasteless?"'; `
       New = "[name=`"Emperor`"]   Who you callin' 'tasteless?'" }
)

for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $cell = $ws.Cells.Item($startRow + $r, $startCol + $c)
        $text = $cell.Value2
        if ($null -eq $text) { continue }
        if ($text -isnot [string]) { continue }

        foreach ($rep in $replacements) {
            if ($text.Contains($rep.Old)) {
                $text = $text.Replace($rep.Old, $rep.New)
            }
        }

        if ($text -ne $cell.Value2) {
            $cell.Value = $text
        }
    }
}
